# Nasar - 24th Nov
# Adds three new rows (97-99) describing the "Biometrics disabled" pop-up
# elements to the iOS object-repository sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS")

# Column B (element values) was populated first for the three new rows...
$ws.Range("B97").Value = '//XCUIElementTypeStaticText[@name="Biometrics disabled"]'
$ws.Range("B98").Value = '//XCUIElementTypeStaticText[@name="Please go to device setting to enable."]'
$ws.Range("B99").Value = "ok"

# ...then column A (element keys), written in this order: 99, 97, 98
$ws.Range("A99").Value = "BTN_OK_IN_POP_UP"
$ws.Range("A97").Value = "BIOMETRICS_DISABLED_LABEL_TXT"
$ws.Range("A98").Value = "PLEASE_GO_TO_DEVICE_SETTINGS_BIO_TXT"

# Finally column E (value type)
$ws.Range("E97").Value = "device-xpath"
$ws.Range("E98").Value = "device-xpath"
$ws.Range("E99").Value = "device-accessibilityid"

# Update the view to match the authored sheet state
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("G59").Select()
